$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells retain their original text (string) storage type,
# matching the source workbook where Price/Volume columns are stored as text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.723.16'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.648.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.53'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0629'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.34'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.56'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.691.42'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.20'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.534'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.23'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.809.99'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0756'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.67'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.40'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.55'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +13.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.24'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.09'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.93'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0523'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.44'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.68%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.62%  '

# Rows 34 and 35 swap places (Maker <-> LidoDAOToken) with updated values
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.27%  '

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.300.50'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.99%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.09%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.15%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.816'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.66%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.06%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.790.27'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.69'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.14'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +11.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.45%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0978'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.408'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.56%  '
